$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data
$lastRow = $ws.Cells.Item(1, 1).End(4).Row

# Mapping for the correct_ans column (L) abbreviation -> full word
$ansMap = @{ "b" = "center"; "y" = "left"; "r" = "right" }

for ($r = 2; $r -le $lastRow; $r++) {

    # Update column L (correct_ans): map abbreviated codes to full words
    $lCell = $ws.Cells.Item($r, 12)
    $lVal = $lCell.Value()
    if ($ansMap.ContainsKey($lVal)) {
        $lCell.Value = $ansMap[$lVal]
    }

    # Update columns A-D (promptFile, correctFile, dist_01File, dist_02File):
    # rename the "face" image category to "book" (folder prefix + filename)
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val.ToString().StartsWith("face//")) {
            $cell.Value = $val.ToString().Replace("face", "book")
        }
    }
}
